# Append the "2021年" row (row 12) to Sheet1, mirroring the layout/formatting
# already used by the preceding data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A12: copy row-11's label cell (value + style) so the new year label keeps
# the same bold/bordered/centered formatting, then overwrite its text.
$ws.Range("A11").Copy($ws.Range("A12"))
$ws.Range("A12").Value = "2021年"

$ws.Range("B12").Value = 2153

# C12 is blank in the source data (stored as an empty text cell, same as the
# equivalent cells in earlier rows). A plain "" assignment would leave the
# cell completely empty (no cell element at all), so force text-typing with
# a leading apostrophe and then strip the resulting quote-prefix formatting.
$ws.Range("C12").Value = "'"
$ws.Range("C12").ClearFormats()

$ws.Range("D12").Value = 3611
$ws.Range("E12").Value = 6827
$ws.Range("F12").Value = 9239
$ws.Range("G12").Value = 23163
$ws.Range("H12").Value = 15391
$ws.Range("I12").Value = 4000
$ws.Range("J12").Value = 1567
$ws.Range("K12").Value = 3545

$ws.Range("L12").Value = "'"
$ws.Range("L12").ClearFormats()

$ws.Range("M12").Value = 142842
$ws.Range("N12").Value = 41159
$ws.Range("O12").Value = 9188
$ws.Range("P12").Value = 1985
$ws.Range("Q12").Value = 2586

$ws.Range("R12").Value = "'"
$ws.Range("R12").ClearFormats()
$ws.Range("S12").Value = "'"
$ws.Range("S12").ClearFormats()

$ws.Range("T12").Value = 743
$ws.Range("U12").Value = 197462
